# childMetrics.xlsx - "about to put in conclusion and discussion edits"
#
# Changes applied:
#  1) Move the active selection from A2:K6 (active cell A2) to a single
#     cell L14 (i.e. the user clicked away from the table, out past the
#     last column, getting ready to type conclusion/discussion notes).
#  2) Re-assert the existing header / label formatting on the table so the
#     workbook's style table is rewritten (bold font + border, with the
#     appropriate alignment) - this mirrors the cosmetic xf bookkeeping
#     seen in the saved file without altering how any cell actually looks.
#  3) Correct the TNR value for the "Log Reg" row (J6) from 0.8857 to 0.75.
#  4) Nudge the saved window position/size, matching the slightly
#     resized/repositioned Excel window in the new file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Selection -----------------------------------------------------
[void]$ws.Range("L14").Select()

# --- 2) Re-apply formatting (visually a no-op, keeps the style table in
#        sync with the rest of the workbook) ---------------------------
$headerRange = $ws.Range("B2:E2,H2:K2")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter

$centerRange = $ws.Range("F2:G2")
$centerRange.Font.Bold = $true
$centerRange.HorizontalAlignment = -4108   # xlCenter
$centerRange.VerticalAlignment = -4108     # xlCenter

$labelRange = $ws.Range("A3:K3,A4,A5,A6")
$labelRange.Font.Bold = $true

# --- 3) Fix the TNR figure for the Log Reg row -------------------------
$ws.Range("J6").Value = 0.75

# --- 4) Match the window geometry recorded in the workbook view --------
$win = $excel.ActiveWindow
$win.Left = 24
$win.Top = 0
$win.Width = 1256
$win.Height = 709
